# Update "想去人数" (want-to-go count) values in column F
# for sheets "展览" (Exhibition) and "全部类型" (All Types).
# Both sheets share the same set of events (except one extra
# performance row in "全部类型"), so the same F-column updates
# apply to both, with row references shifted by 1 starting at
# row 32 in "全部类型" because of that extra row.

$wb = $excel.ActiveWorkbook

# Row -> (old, new) value map for the "展览" sheet
$sheet1Updates = @{
    2  = 117
    3  = 308
    4  = 61
    5  = 755
    6  = 73
    7  = 2116
    8  = 276
    10 = 4686
    11 = 5
    12 = 37
    14 = 216
    15 = 20
    16 = 155
    17 = 31
    18 = 21
    19 = 100
    20 = 3592
    21 = 173
    22 = 579
    25 = 93
    26 = 106
    27 = 14
    28 = 8
    32 = 779
    33 = 2222
}

# Row -> value map for the "全部类型" sheet (rows 32/33 shifted to 33/34)
$sheet4Updates = @{
    2  = 117
    3  = 308
    4  = 61
    5  = 755
    6  = 73
    7  = 2116
    8  = 276
    10 = 4686
    11 = 5
    12 = 37
    14 = 216
    15 = 20
    16 = 155
    17 = 31
    18 = 21
    19 = 100
    20 = 3592
    21 = 173
    22 = 579
    25 = 93
    26 = 106
    27 = 14
    28 = 8
    33 = 779
    34 = 2222
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
